# Applies the "第九周 周二" (week 9, Tuesday) plan block to the sheet,
# mirroring the structure/styling of the preceding week block (rows 142-150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting of the previous week's block (header + table + summary)
#    down onto the new block of rows so fonts/fills/borders/alignment match exactly,
#    without creating duplicate style entries.
$ws.Range("A142:D150").Copy()
$ws.Range("A152:D160").PasteSpecial(-4122)

# 2) Fill in the new block's text content.

# Date / header row (merged A152:D152)
$ws.Range("A152").Value = "日期：2017.10.23 第九周 周二"

# Column headers row
$ws.Range("A153").Value = "人员"
$ws.Range("B153").Value = "计划任务"
$ws.Range("C153").Value = "完成情况"
$ws.Range("D153").Value = "备注"

# Per-person task rows
$ws.Range("A154").Value = "李杰"
$ws.Range("B154").Value = "开发后台“用户管理”模块中的控制层(controller)"

$ws.Range("A155").Value = "周振朋"
$ws.Range("B155").Value = "开发前端首页的“分类”模块"

$ws.Range("A156").Value = "禤锦辉"
$ws.Range("B156").Value = "设计图标，添加订单表和商品类别表数据"

$ws.Range("A157").Value = "柯新钿"
$ws.Range("B157").Value = "开发前端“我的”模块"

$ws.Range("A158").Value = "冯文雄"
$ws.Range("B158").Value = "开发后台“宝贝管理”模块中的控制层(controller)"

$ws.Range("A159").Value = "阿卜力孜"
$ws.Range("B159").Value = "添加数据库的宝贝留言表和求购留言表的数据"

# Summary row (merged A160:D160)
$ws.Range("A160").Value = "总结："

# 3) Row heights for the new block.
$ws.Rows.Item(152).RowHeight = 22.5
$ws.Rows.Item(153).RowHeight = 22.5
$ws.Rows.Item(154).RowHeight = 45
$ws.Rows.Item(155).RowHeight = 45
$ws.Rows.Item(156).RowHeight = 45
$ws.Rows.Item(157).RowHeight = 22.5
$ws.Rows.Item(158).RowHeight = 45
$ws.Rows.Item(159).RowHeight = 45
$ws.Rows.Item(160).RowHeight = 22.5

# 4) Merge the header/date row and the summary row, same as other week blocks.
$ws.Range("A152:D152").Merge()
$ws.Range("A160:D160").Merge()

# 5) Match the saved selection/active cell from the source workbook.
$ws.Range("B154").Select()
